# suppression details calendrier + ajout doc sante
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calendrier")

# Row 39 - "Prix Geko Bike Road à Bernwiller": remove the "Ouvert aux FFC" note (col G)
$ws.Range("G39").Value = ""

# Row 41 - "Grand Prix d'Ampfersbach": shorten the cancellation note (col G)
$ws.Range("G41").Value = "Annulé"

# Row 42 - "20e Grand Prix Gestimmo à Magstatt le Bas": shorten the cancellation note (col G)
$ws.Range("G42").Value = "Annulé"

# Row 59 - "Gentlemen de Nommay": drop the "- à confirmer" suffix (col B)
$ws.Range("B59").Value = "Gentlemen de Nommay (épreuve FFC ouverte aux FSGT)"

# Row 59 - add the FFC results/health-pass document link (col F)
$ws.Range("F59").Value = "https://www.ffc-bfc.fr/calendrier-resultats/gentlemen-de-nommay-2"

# Leave the cursor where the author ended up after the edit
$ws.Range("B60").Select()
